# Applies the TestCases.xlsx edit:
#  - TestCases sheet: delete row 59 (PESmokeTc044_3_verifyAccountAssociationImpact),
#    shifting all rows below up by one.
#  - TestCases!E38 (Priority) changes from "High" to "Low".
#  - TestCases!D113 (Status, after the row shift) changes to "Pass".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# Update priority on row 38 before the shift happens (row 38 is above the
# deleted row, so its row number never changes).
$ws.Range("E38").Value = "Low"

# Delete the entire row for PESmokeTc044_3_verifyAccountAssociationImpact.
# This shifts every row below it up by one (row 60 becomes 59, ... row 114
# becomes 113).
$ws.Rows(59).Delete()

# The last test case row (originally row 114, now row 113) has its Status
# updated to "Pass".
$ws.Range("D113").Value = "Pass"
